$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 'MSG: None

MSG: The decision has been recorded to acquire the rights for "Barbie" to be shown on Friday.
'
$ws.Range("C3").Value = 'MSG: None

MSG: The decision has been made to acquire the rights for "Oppenheimer."
'
$ws.Range("C4").Value = 'MSG: None

MSG: The rights to both movies have been acquired.
'
$ws.Range("D4").Value = 'both_movies, '
$ws.Range("C5").Value = 'MSG: None

MSG: The conversation has concluded without a decision on which movie to show on Friday.
'
$ws.Range("D5").Value = 'no_decision, '
$ws.Range("C6").Value = 'MSG: None

MSG: The decision has been made to acquire the rights to show "Barbie."
'
$ws.Range("C7").Value = 'MSG: None

MSG: The rights to both movies have been successfully acquired for the upcoming screenings.
'
$ws.Range("D7").Value = 'both_movies, '
$ws.Range("C8").Value = 'MSG: None

MSG: The decision has been recorded, and no movie will be acquired as the committee did not reach a definitive agreement.
'
$ws.Range("D8").Value = 'no_decision, '
$ws.Range("C9").Value = 'MSG: None

MSG: The decision was made that no movie would be shown on Friday.
'
$ws.Range("D9").Value = 'no_decision, '
$ws.Range("C10").Value = 'MSG: None

MSG: The decision about what movie to show on Friday was not reached.
'
$ws.Range("D10").Value = 'no_decision, '
$ws.Range("C11").Value = 'MSG: None

MSG: The decision about what movie to show on Friday has not been made.
'
$ws.Range("D11").Value = 'no_decision, '
$ws.Range("C12").Value = 'MSG: None

MSG: The committee did not reach a decision regarding the movie to show on Friday.
'
$ws.Range("D12").Value = 'no_decision, '
$ws.Range("C13").Value = 'MSG: None

MSG: The decision from the committee has been recorded and "Barbie" was selected as the movie to acquire for Friday''s showing.
'
$ws.Range("C14").Value = 'MSG: None

MSG: The function call for no decision has been made, indicating that the committee did not reach a consensus on which movie to show on Friday.
'
$ws.Range("D14").Value = 'no_decision, '
$ws.Range("C15").Value = 'MSG: None

MSG: The decision has been recorded, and "Barbie" will be acquired for showing on Friday.
'
$ws.Range("C16").Value = 'MSG: None

MSG: The decision has been recorded to acquire the rights for "Barbie."
'
$ws.Range("C17").Value = 'MSG: None

MSG: The decision-making process did not yield an agreement on which movie to show on Friday, so the outcome is that no decision was made.
'
$ws.Range("D17").Value = 'no_decision, '
$ws.Range("C18").Value = 'MSG: None

MSG: The decision about which movie to show on Friday could not be made.
'
$ws.Range("D18").Value = 'no_decision, '
$ws.Range("C19").Value = 'MSG: None

MSG: The decision to acquire the rights for "Barbie" has been successfully recorded. If you have any further questions or need assistance, feel free to ask!
'
$ws.Range("C20").Value = 'MSG: None

MSG: The decision-making process ended without a definitive choice for Friday''s movie.
'
$ws.Range("D20").Value = 'no_decision, '
$ws.Range("C21").Value = 'MSG: None

MSG: It appears that there was no decision made regarding which movie to show on Friday. Therefore, I will proceed with indicating that no decision has been reached.
'
$ws.Range("D21").Value = 'no_decision, '
$ws.Range("C22").Value = 'MSG: None

MSG: The decision about which movie to show on Friday has not been made, so I will call the no_decision function.
'
$ws.Range("D22").Value = 'no_decision, '
$ws.Range("C23").Value = 'MSG: None

MSG: The decision has been recorded, and no movie was selected for the meeting on Friday.
'
$ws.Range("D23").Value = 'no_decision, '
$ws.Range("C24").Value = 'MSG: None

MSG: The decision to acquire the rights for both movies has been recorded.
'
$ws.Range("C25").Value = 'MSG: None

MSG: The decision has been recorded as "no decision."
'
$ws.Range("D25").Value = 'no_decision, '
$ws.Range("C26").Value = 'MSG: None

MSG: The decision to show a movie on Friday was ultimately not made, leading to no acquisition of movie rights.
'
$ws.Range("D26").Value = 'no_decision, '
$ws.Range("C27").Value = 'MSG: None

MSG: The decision to show "Barbie" has been confirmed.
'
$ws.Range("C28").Value = 'MSG: None

MSG: The decision has been recorded: no movie was selected for showing on Friday.
'
$ws.Range("D28").Value = 'no_decision, '
$ws.Range("C29").Value = 'MSG: None

MSG: The decision has been recorded to acquire the rights for "Barbie."
'
$ws.Range("C30").Value = 'MSG: None

MSG: The decision has been recorded as "no decision" regarding the movie to be shown on Friday.
'
$ws.Range("D30").Value = 'no_decision, '
$ws.Range("C31").Value = 'MSG: None

MSG: The decision process ended without a plan about which movie to show on Friday, resulting in no agreement.
'
$ws.Range("D31").Value = 'no_decision, '
$ws.Range("C32").Value = 'MSG: None

MSG: The decision has been made to acquire the rights for "Barbie."
'
$ws.Range("C33").Value = 'MSG: None

MSG: The decision has been recorded as no agreement was made on which movie to show on Friday.
'
$ws.Range("D33").Value = 'no_decision, '
$ws.Range("C34").Value = 'MSG: None

MSG: The decision has been recorded as no movie was selected for Friday.
'
$ws.Range("D34").Value = 'no_decision, '
$ws.Range("C35").Value = 'MSG: None

MSG: The decision has been recorded as no consensus on the selection of a movie for Friday.
'
$ws.Range("D35").Value = 'no_decision, '
$ws.Range("C36").Value = 'MSG: None

MSG: The decision regarding which movie to show on Friday has not been made, so the function indicating no decision has been executed.
'
$ws.Range("D36").Value = 'no_decision, '
$ws.Range("C37").Value = 'MSG: None

MSG: The decision about which movie to show on Friday has resulted in no agreement.
'
$ws.Range("D37").Value = 'no_decision, '
$ws.Range("C38").Value = 'MSG: None

MSG: The committee did not reach a decision about which movie to show on Friday.
'
$ws.Range("D38").Value = 'no_decision, '
$ws.Range("C39").Value = 'MSG: None

MSG: The decision has been recorded, and the rights to "Barbie" will be acquired.
'
$ws.Range("C40").Value = 'MSG: None

MSG: The decision has been recorded to acquire the rights for "Oppenheimer."
'
$ws.Range("C41").Value = 'MSG: None

MSG: No decision was made about which movie to show on Friday.
'
$ws.Range("D41").Value = 'no_decision, '
$ws.Range("C42").Value = 'MSG: None

MSG: The decision has been recorded: no movie was selected for Friday.
'
$ws.Range("D42").Value = 'no_decision, '
$ws.Range("C43").Value = 'MSG: None

MSG: The rights for both movies have been successfully acquired.
'
$ws.Range("D43").Value = 'both_movies, '
$ws.Range("C44").Value = 'MSG: None

MSG: The decision about the movie for Friday was ultimately not reached, resulting in no agreement.
'
$ws.Range("D44").Value = 'no_decision, '
$ws.Range("C45").Value = 'MSG: None

MSG: The decision to select "Barbie" for the event has been recorded.
'
$ws.Range("C46").Value = 'MSG: None

MSG: The decision regarding the movie to show on Friday has not been made.
'
$ws.Range("D46").Value = 'no_decision, '
$ws.Range("C47").Value = 'MSG: None

MSG: The decision has been recorded, and "Barbie" will be the movie shown on Friday.
'
$ws.Range("C48").Value = 'MSG: None

MSG: The decision has been recorded, and there was no agreement on a movie to show on Friday.
'
$ws.Range("D48").Value = 'no_decision, '
$ws.Range("C49").Value = 'MSG: None

MSG: The decision has been recorded as no decision. If further discussion occurs, please let me know!
'
$ws.Range("D49").Value = 'no_decision, '
$ws.Range("C50").Value = 'MSG: None

MSG: The rights to both movies have been acquired.
'
$ws.Range("D50").Value = 'both_movies, '
$ws.Range("C51").Value = 'MSG: None

MSG: The decision has been recorded as "no decision" regarding the movie to be shown on Friday.
'
$ws.Range("D51").Value = 'no_decision, '
$ws.Range("C52").Value = 'MSG: None

MSG: The decision has been recorded to acquire the rights for the movie "Barbie."
'
$ws.Range("C53").Value = 'MSG: None

MSG: The decision has been recorded as a "no decision" regarding the movie selection for Friday.
'
$ws.Range("D53").Value = 'no_decision, '
$ws.Range("C54").Value = 'MSG: None

MSG: The decision has been recorded, and the rights for "Barbie" will be acquired for the movie to be shown on Friday.
'
$ws.Range("C55").Value = 'MSG: None

MSG: I have successfully recorded the decision to acquire rights for both movies, "Oppenheimer" and "Barbie," to be shown on Friday.
'
$ws.Range("D55").Value = 'both_movies, '
$ws.Range("C56").Value = 'MSG: None

MSG: The decision about which movie to show on Friday was not made.
'
$ws.Range("D56").Value = 'no_decision, '
$ws.Range("C57").Value = 'MSG: None

MSG: The decision has been recorded as no movie being selected.
'
$ws.Range("D57").Value = 'no_decision, '
$ws.Range("C58").Value = 'MSG: None

MSG: No decision about Friday''s movie was made.
'
$ws.Range("D58").Value = 'no_decision, '
